$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.762.33"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "1.806.10"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.557"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.288"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("E10").Value = "  +8.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "2.063.34"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "1.789.35"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "34.760.98"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("E20").Value = "  +9.45%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0536"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "1.449.76"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  +3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0494"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.960.30"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.88%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0127"
$ws.Range("E51").Value = "  +9.74%  "
